$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column E
$ws.Range("E1").Value = "Functional Kalman Python"
$ws.Range("E1").Font.Bold = $true

# New column width to match column B
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Data values for column E (rows 2-11)
$values = @(355.01, 77.25, 91.08, 60.08, 60.08, 66.04, 82.97, 83.21, 94.89, 56.74)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}

# Average formula in E12
$ws.Range("E12").Formula = "=AVERAGE(E2:E11)"

# Update selection to F1, matching the diff
$ws.Range("F1").Select()
